$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table with the latest scraped values.
# For cells whose new text looks like a plain decimal number (e.g. "0.998",
# "12.70"), the cell's NumberFormat is first forced to Text ("@") so Excel
# does not silently reinterpret the string as a Number (which would drop
# meaningful trailing zeros / formatting). Values that already contain
# extra punctuation (thousand separators, subscripts, etc.) stay text
# naturally and do not need this treatment.

$ws.Range("D2").Value = '34.640.50'
$ws.Range("E2").Value = '  +13.12%  '
$ws.Range("D3").Value = '1.846.29'
$ws.Range("E3").Value = '  +10.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.36'
$ws.Range("E5").Value = '  +5.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.580'
$ws.Range("E6").Value = '  +9.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.11'
$ws.Range("E8").Value = '  +8.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.14'
$ws.Range("E9").Value = '  +6.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.292'
$ws.Range("E10").Value = '  +10.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0695'
$ws.Range("E11").Value = '  +7.78%  '
$ws.Range("E12").Value = '  +3.42%  '
$ws.Range("D13").Value = '2.111.11'
$ws.Range("E13").Value = '  +10.09%  '
$ws.Range("D14").Value = '1.852.92'
$ws.Range("E14").Value = '  +10.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.662'
$ws.Range("E15").Value = '  +8.57%  '
$ws.Range("D16").Value = '34.613.48'
$ws.Range("E16").Value = '  +12.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '10.43'
$ws.Range("E17").Value = '  +4.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.36'
$ws.Range("E18").Value = '  +8.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.08'
$ws.Range("E19").Value = '  +7.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '262.35'
$ws.Range("E20").Value = '  +8.16%  '
$ws.Range("D21").Value = '0.0₃0770'
$ws.Range("E21").Value = '  +6.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.75'
$ws.Range("E23").Value = '  +7.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.45'
$ws.Range("E24").Value = '  +4.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("E25").Value = '  +4.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.54'
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.08'
$ws.Range("E27").Value = '  +8.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.120'
$ws.Range("E28").Value = '  +6.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.27'
$ws.Range("E29").Value = '  +8.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.93'
$ws.Range("E31").Value = '  +13.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0532'
$ws.Range("E32").Value = '  +7.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.23'
$ws.Range("E33").Value = '  +7.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.65'
$ws.Range("E34").Value = '  +10.94%  '
$ws.Range("D35").Value = '1.566.03'
$ws.Range("E35").Value = '  +4.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.84'
$ws.Range("E36").Value = '  +3.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.09'
$ws.Range("E37").Value = '  +6.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.652'
$ws.Range("E38").Value = '  +9.09%  '
$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '86.86'
$ws.Range("E39").Value = '  +3.03%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0192'
$ws.Range("E40").Value = '  +7.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.85'
$ws.Range("E41").Value = '  +6.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.932'
$ws.Range("E42").Value = '  +11.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.34'
$ws.Range("E43").Value = '  +2.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.17'
$ws.Range("E44").Value = '  +7.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.14'
$ws.Range("E45").Value = '  +175.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0528'
$ws.Range("E46").Value = '  +5.33%  '
$ws.Range("E47").Value = '  +5.99%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.70'
$ws.Range("E48").Value = '  +24.38%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.001.05'
$ws.Range("E49").Value = '  +10.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.89'
$ws.Range("E50").Value = '  +6.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("E51").Value = '  -0.24%  '
